# Raw and Clean Data from SSA for October 16-18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows before old row 134 (old 134-136 shift down to 137-139) ---
$ws.Rows("134:136").Insert()

# Row 134: 2020-10-10 (raw data row, default font, date number format)
$ws.Range("A134").Value = 44114
$ws.Range("A134").NumberFormat = "mm-dd-yy"
$ws.Range("B134:F134").Font.Size = 12
$ws.Range("B134:F134").WrapText = $true
$ws.Range("B134").Value = 814328
$ws.Range("C134").Value = 964812
$ws.Range("D134").Value = 305487
$ws.Range("E134").Value = 83642
$ws.Range("F134").Value = 23.404205676336808
$ws.Rows(134).RowHeight = 16

# Row 135: 2020-10-11
$ws.Range("A135").Value = 44115
$ws.Range("A135").NumberFormat = "mm-dd-yy"
$ws.Range("B135:F135").Font.Size = 12
$ws.Range("B135:F135").WrapText = $true
$ws.Range("B135").Value = 817503
$ws.Range("C135").Value = 969859
$ws.Range("D135").Value = 301579
$ws.Range("E135").Value = 83781
$ws.Range("F135").Value = 23.367865316702201
$ws.Rows(135).RowHeight = 16

# Row 136: 2020-10-12 (date cell rendered in black font explicitly)
$ws.Range("A136").Value = 44116
$ws.Range("A136").Font.Color = 0
$ws.Range("A136").NumberFormat = "mm-dd-yy"
$ws.Range("B136:F136").Font.Size = 12
$ws.Range("B136:F136").WrapText = $true
$ws.Range("B136").Value = 821045
$ws.Range("C136").Value = 975299
$ws.Range("D136").Value = 297064
$ws.Range("E136").Value = 83945
$ws.Range("F136").Value = 23.32965915388316
$ws.Rows(136).RowHeight = 16

# --- New rows appended after the (now shifted) old tail rows 137-139 ---
# Row 140: 2020-10-16
$ws.Range("A140").Value = 44120
$ws.Range("A140").Font.Size = 12
$ws.Range("A140").WrapText = $true
$ws.Range("A140").NumberFormat = "mm-dd-yy"
$ws.Range("B140:F140").Font.Size = 12
$ws.Range("B140:F140").WrapText = $true
$ws.Range("B140").Value = 841661
$ws.Range("C140").Value = 1004800
$ws.Range("D140").Value = 313678
$ws.Range("E140").Value = 85704
$ws.Range("F140").Value = 23.229780160896134
$ws.Rows(140).RowHeight = 16

# Row 141: 2020-10-17
$ws.Range("A141").Value = 44121
$ws.Range("A141").Font.Size = 12
$ws.Range("A141").WrapText = $true
$ws.Range("A141").NumberFormat = "mm-dd-yy"
$ws.Range("B141:F141").Font.Size = 12
$ws.Range("B141:F141").WrapText = $true
$ws.Range("B141").Value = 847108
$ws.Range("C141").Value = 1013186
$ws.Range("D141").Value = 316228
$ws.Range("E141").Value = 86059
$ws.Range("F141").Value = 23.186063642416315
$ws.Rows(141).RowHeight = 16

# Row 142: 2020-10-18 (plain text date label like the other SSA-report rows)
$ws.Range("A142").Value = "2020-10-18"
$ws.Range("B142").Value = 851227
$ws.Range("C142").Value = 1019821
$ws.Range("D142").Value = 312224
$ws.Range("E142").Value = 86167
$ws.Range("F142").Value = 23.32

# --- View state: scroll/selection to match the authored workbook state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 104
$win.ScrollColumn = 1
$ws.Range("L133").Select()
